# Documentation correction for Umbraco >7.00
#
# This script applies the textual corrections described in the commit:
#  - The key that needs to be added is now called "appSettings" and can live
#    either in Web.config or in appSettings.config, depending on the Umbraco
#    setup, instead of always living in "Config/appSettings.config".
#  - Every other mention of the old "Config/appSettings.config file" wording
#    is simplified to just "appSettings".

$d = $word.ActiveDocument

# 1) First / most detailed mention, right after explaining the
#    "umbracoNaviHide" alias. Expand "the Config/appSettings.config file:"
#    into the longer explanation that references both Web.config and
#    appSettings.config.
$d.Content.Find.Execute( `
    "needs to be added to the Config/appSettings.config file:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "needs to be added to appSettings (in Web.config or appSettings.config, depending on your Umbraco setup):", `
    2) | Out-Null

# 2) "...can be added to the web site's  Config/appSettings.config file)."
$d.Content.Find.Execute( `
    "Config/appSettings.config file).", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "appSettings).", `
    2) | Out-Null

# 3) "...by entering the following to the Config/appSettings.config file:"
$d.Content.Find.Execute( `
    "to the Config/appSettings.config file:", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "to appSettings:", `
    2) | Out-Null

# 4) "...can be added to the web site's Config/appSettings.config file." (ends a paragraph)
$d.Content.Find.Execute( `
    "web site's Config/appSettings.config file.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "web site's appSettings.", `
    2) | Out-Null

# 5) "...needs to be entered into the web site's Config/appSettings.config file. Then,"
$d.Content.Find.Execute( `
    "web site's Config/appSettings.config file. Then,", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "web site's appSettings. Then,", `
    2) | Out-Null
